{"js": "// Fix the grammatical agreement (\"mani\u00e8re\" + adjective) in each of the\n// five cadence-heading lines of the document (\"constant(e)\",\n// \"journali\u00e8re(s)\", \"hebdomadaire(s)\", \"mensuel(le)\", \"annuel(le)\").\nconst replacements = [\n  [\"De mani\u00e8re constant:\", \"De mani\u00e8re constante:\"],\n  [\"De mani\u00e8res journali\u00e8res:\", \"De mani\u00e8re journali\u00e8re:\"],\n  [\"De mani\u00e8res hebdomadaires:\", \"De mani\u00e8re hebdomadaire:\"],\n  [\"De mani\u00e8res mensuel:\", \"De mani\u00e8re mensuelle:\"],\n  [\"De mani\u00e8res annuel:\", \"De mani\u00e8re annuelle:\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items/text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: \"${oldText}\"`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Fix the grammatical agreement (\"mani\u00e8re\" + adjective) in each of the\n# five cadence-heading lines of the document (\"constant(e)\",\n# \"journali\u00e8re(s)\", \"hebdomadaire(s)\", \"mensuel(le)\", \"annuel(le)\").\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"De mani\u00e8re constant:\", \"De mani\u00e8re constante:\"),\n    @(\"De mani\u00e8res journali\u00e8res:\", \"De mani\u00e8re journali\u00e8re:\"),\n    @(\"De mani\u00e8res hebdomadaires:\", \"De mani\u00e8re hebdomadaire:\"),\n    @(\"De mani\u00e8res mensuel:\", \"De mani\u00e8re mensuelle:\"),\n    @(\"De mani\u00e8res annuel:\", \"De mani\u00e8re annuelle:\")\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}\n"}
